$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("legislators_data_dictionary")

# Insert a new row above the current row 18 ("years_active"), shifting
# everything down by one, and fill in the new "wiki_url" field row.
$ws.Rows.Item(18).Insert()

$ws.Range("B18").Value = "TEXT "
$ws.Range("A18").Value = "wiki_url (UNIQUE)"
$ws.Range("D18").Value = "Wikipedia URL for the legislator, used to uniquely identify legislators across the country"
$ws.Range("E18").Value = "https://en.wikipedia.org/wiki/Jeremy_Harper_(politician)"

# Excel auto-restyles a cell whose text looks like a URL; put the
# quote-prefixed "plain text" style back so it matches the rest of the
# Examples column (same style as the "gender" example directly above it).
$ws.Range("E18").Style = "Normal"
$ws.Cells.Item(18, 5).Font.Name = $ws.Cells.Item(17, 5).Font.Name

$ws.Rows.Item(18).RowHeight = 30

$wb.Save()
